$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Bump individual stack-trace line numbers (M2DocEvaluator / M2DocUtils / AbstractTemplatesTestSuite)
Replace-Text "caseQuery(M2DocEvaluator.java:559)" "caseQuery(M2DocEvaluator.java:586)"

# "doSwitch(M2DocEvaluator.java:1216)" -> "...:1239)" occurs 3 times with the same replacement.
Replace-Text "M2DocEvaluator.java:1216)" "M2DocEvaluator.java:1239)"

Replace-Text "caseBlock(M2DocEvaluator.java:1425)" "caseBlock(M2DocEvaluator.java:1464)"
Replace-Text "caseDocumentTemplate(M2DocEvaluator.java:287)" "caseDocumentTemplate(M2DocEvaluator.java:296)"
Replace-Text "generate(M2DocEvaluator.java:276)" "generate(M2DocEvaluator.java:281)"
Replace-Text "M2DocUtils.generate(M2DocUtils.java:694)" "M2DocUtils.generate(M2DocUtils.java:805)"
Replace-Text "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)" "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)"
Replace-Text "generation(AbstractTemplatesTestSuite.java:389)" "generation(AbstractTemplatesTestSuite.java:420)"

# Insert a new stack frame line right before the unique
# "RunAfters.evaluate(RunAfters.java:27)" that directly follows
# "ParentRunner$2.evaluate(ParentRunner.java:268)" (there are several
# occurrences of that ParentRunner line in the trace, but only one is
# immediately followed by RunAfters.evaluate rather than ParentRunner.run).
$tab = [char]9
$nl = [char]10
$old = "ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl + $tab + "at org.junit.internal.runners.statements.RunAfters"
$new = "ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl + $tab + "at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)" + $nl + $tab + "at org.junit.internal.runners.statements.RunAfters"
Replace-Text $old $new
